# Daily attendance processing - 2026-01-02 11:31:35
# Normalizes the "Recorded By" (column G) cell text by reordering the
# comma-separated list of recorders for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value is exactly "dnasr281@gmail.com, System"
# and must become "System, dnasr281@gmail.com"
$rowsSwap1 = @(3,6,10,12,13,14,15,18,19,20,21,22,24,26,29,32,36,38,39,40,41,44,45,46,47,48,50,52,55,58,62,64,65,66,67,70,71,72,73,74,76,78,83,84,85,86,90,92,99,101,109,110,111,112,116,118,125,127,135,136,137,138,142,144,151,153)

foreach ($r in $rowsSwap1) {
    $cell = $ws.Range("G$r")
    if ($cell.Value() -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# Rows whose "Recorded By" value is exactly "system, backup@backdoor.com, System"
# and must become "backup@backdoor.com, System, system"
$rowsSwap2 = @(2,28,54)

foreach ($r in $rowsSwap2) {
    $cell = $ws.Range("G$r")
    if ($cell.Value() -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
